# Update crypto price/volume figures per the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.935.89'
$ws.Range('E2').Value = '  -4.66%  '
$ws.Range('D3').Value = '2.223.96'
$ws.Range('E3').Value = '  -5.68%  '
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '317.69'
$c.ClearFormats()
$ws.Range('E5').Value = '  +2.10%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '100.42'
$c.ClearFormats()
$ws.Range('E6').Value = '  -6.60%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.592'
$c.ClearFormats()
$ws.Range('E7').Value = '  -5.96%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -6.78%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '37.41'
$c.ClearFormats()
$ws.Range('E10').Value = '  -8.12%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '54.03'
$c.ClearFormats()
$ws.Range('E11').Value = '  -2.98%  '
$ws.Range('E12').Value = '  -8.95%  '
$ws.Range('E13').Value = '  -7.18%  '
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('E15').Value = '  -10.81%  '
$ws.Range('D16').Value = '2.562.54'
$ws.Range('E16').Value = '  -5.60%  '
$ws.Range('E17').Value = '  -5.46%  '
$ws.Range('D18').Value = '2.230.71'
$ws.Range('E18').Value = '  -5.52%  '
$ws.Range('D19').Value = '42.852.28'
$ws.Range('E19').Value = '  -4.75%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '15.07'
$c.ClearFormats()
$ws.Range('E20').Value = '  +5.43%  '
$ws.Range('E21').Value = '  -8.43%  '
$ws.Range('E22').Value = '  -10.00%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '65.60'
$c.ClearFormats()
$ws.Range('E23').Value = '  -9.81%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '3.17'
$c.ClearFormats()
$ws.Range('E24').Value = '  -9.57%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '236.84'
$c.ClearFormats()
$ws.Range('E25').Value = '  -8.16%  '
$ws.Range('E26').Value = '  -6.75%  '
$ws.Range('E27').Value = '  -0.35%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '10.10'
$c.ClearFormats()
$ws.Range('E28').Value = '  -8.50%  '
$ws.Range('E29').Value = '  -4.64%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '6.42'
$c.ClearFormats()
$ws.Range('E30').Value = '  -10.38%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.0911'
$c.ClearFormats()
$ws.Range('E31').Value = '  -5.73%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '20.51'
$c.ClearFormats()
$ws.Range('E32').Value = '  -7.83%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '34.26'
$c.ClearFormats()
$ws.Range('E33').Value = '  -7.74%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '157.26'
$c.ClearFormats()
$ws.Range('E34').Value = '  -6.18%  '
$ws.Range('E35').Value = '  -6.53%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '3.19'
$c.ClearFormats()
$ws.Range('E36').Value = '  +10.32%  '
$ws.Range('E37').Value = '  +13.61%  '
$ws.Range('E38').Value = '  -5.58%  '
$ws.Range('E39').Value = '  +0.95%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '4.48'
$c.ClearFormats()
$ws.Range('E40').Value = '  -3.91%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.106'
$c.ClearFormats()
$ws.Range('E41').Value = '  -8.37%  '
$ws.Range('E42').Value = '  -6.81%  '
$ws.Range('D43').Value = '1.951.71'
$ws.Range('E43').Value = '  +3.57%  '
$ws.Range('E44').Value = '  -0.11%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '12.50'
$c.ClearFormats()
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('E46').Value = '  -11.22%  '
$ws.Range('E47').Value = '  -8.61%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '5.39'
$c.ClearFormats()
$ws.Range('E48').Value = '  -4.07%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '76.57'
$c.ClearFormats()
$ws.Range('E49').Value = '  -5.66%  '
$ws.Range('E50').Value = '  -12.24%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.877'
$c.ClearFormats()
$ws.Range('E51').Value = '  +19.76%  '
